$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in new training results: mask_rcnn_R101_FPN_1x_coco2014_train_valminusmini, iter 90k (row 14 = Box) ---
$ws.Range("E14").Value = 0.3068
$ws.Range("F14").Value = 0.5246
$ws.Range("G14").Value = 0.3204
$ws.Range("H14").Value = 0.1648
$ws.Range("I14").Value = 0.3402
$ws.Range("J14").Value = 0.4068

# --- Fill in new testing data coco2014_minival results (row 15 = Mask) ---
$ws.Range("E15").Value = 0.2764
$ws.Range("F15").Value = 0.484
$ws.Range("G15").Value = 0.2809
$ws.Range("H15").Value = 0.1097
$ws.Range("I15").Value = 0.3024
$ws.Range("J15").Value = 0.4202

# Row 15's number cells lose their inherited border/alignment formatting
# (match plain font-only style used elsewhere in the sheet).
$dataRow15 = $ws.Range("E15:J15")
$dataRow15.Borders.LineStyle = -4142
$dataRow15.HorizontalAlignment = 1
$dataRow15.VerticalAlignment = -4107

# --- Re-seat the row 12-17 merges (unmerge/re-merge) ---
$ws.Range("B14:B15").UnMerge()
$ws.Range("B14:B15").Merge()
$ws.Range("C14:C15").UnMerge()
$ws.Range("C14:C15").Merge()
$ws.Range("B16:B17").UnMerge()
$ws.Range("B16:B17").Merge()
$ws.Range("C16:C17").UnMerge()
$ws.Range("C16:C17").Merge()
$ws.Range("C12:C13").UnMerge()
$ws.Range("C12:C13").Merge()

# --- Move the view's active selection ---
$ws.Range("G22").Select()
